$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A60").NumberFormat = "@"
$ws.Range("A60").Value = "01-07-2021"
$ws.Range("A60").ClearFormats()
$ws.Range("B60").Value = -0.38
$ws.Range("C60").Value = -0.06
$ws.Range("D60").Value = 0.16
